# Update posts.xlsx after post
#
# The post that used to live in row 534 ("「選挙」إنتخابات ...") was removed.
# Deleting the entire row shifts every following row up by one, which is
# exactly what the target diff shows (old row N -> new row N-1 for N >= 535),
# and the sheet's used range shrinks from A1:C597 to A1:C596.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(534).Delete()

$wb.Save()
